$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like updates (coin names and coinranking.com links) -- plain text values
$textUpdates = @{
    'B3' = 'OKB'
    'C3' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'B4' = 'HuobiToken'
    'C4' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'B5' = 'Cronos'
    'C5' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'B6' = 'FTXToken'
    'C6' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B12' = 'MCDex'
    'C12' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'B13' = 'MandalaExchangeToken'
    'C13' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B17' = 'CoinExToken'
    'C17' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'B18' = 'TigerCash'
    'C18' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
}

# Numeric-looking updates (price / volume %) -- must stay text, matching the
# original inline-string cells, so we force the cell format to Text before
# assigning, then restore the default (Normal) style so no stray formatting
# is left behind.
$numericLikeUpdates = @{
    'D2' = '328.29'
    'E2' = '0.98%'
    'D3' = '44.04'
    'E3' = '-1.31%'
    'D4' = '5.512'
    'E4' = '0.34%'
    'D5' = '0.08016'
    'E5' = '-0.33%'
    'D6' = '1.991'
    'E6' = '5.68%'
    'D7' = '4.334'
    'E7' = '0.91%'
    'D8' = '2.613'
    'E8' = '-3.68%'
    'D9' = '0.9494'
    'E9' = '1.32%'
    'D10' = '0.1130'
    'E10' = '-3.66%'
    'D11' = '0.1861'
    'E11' = '-0.24%'
    'D12' = '10.67'
    'E12' = '25.02%'
    'D13' = '0.09807'
    'E13' = '-1.33%'
    'D14' = '0.04704'
    'E14' = '10.55%'
    'D15' = '0.1065'
    'E15' = '0.09%'
    'D16' = '0.001271'
    'E16' = '-0.29%'
    'D17' = '0.04062'
    'E17' = '-4.35%'
    'D18' = '0.005917'
    'E18' = '1.12%'
    'D19' = '3.357'
    'E19' = '-6.54%'
    'D20' = '0.3476'
    'E20' = '-0.25%'
    'E21' = '3.72%'
    'D22' = '0.2545'
    'E22' = '-4.24%'
    'E23' = '1.69%'
    'D24' = '0.004325'
    'E24' = '-3.46%'
    'E25' = '-0.15%'
    'D26' = '0.0003742'
    'E26' = '-6.28%'
    'D38' = '0.02576'
    'E38' = '-2.78%'
    'D39' = '0.05645'
    'E39' = '2.71%'
    'D40' = '0.007541'
    'E40' = '-1.75%'
    'D41' = '0.1397'
    'E41' = '0.35%'
    'D42' = '0.007505'
    'E42' = '1.35%'
    'D43' = '0.002014'
    'E43' = '-1.68%'
    'D44' = '0.008585'
    'E44' = '-1.21%'
    'D45' = '0.00007136'
    'E45' = '0.38%'
    'E46' = '-0.15%'
    'E47' = '55.30%'
    'D48' = '0.003609'
    'E48' = '2.28%'
    'D49' = '0.00002099'
    'E49' = '-0.15%'
    'D50' = '0.0001999'
    'E50' = '-0.15%'
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericLikeUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLikeUpdates[$ref]
    $cell.Style = "Normal"
}
